$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44284
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 500

# Row 3
$ws.Range("D3").Value = 44291
$ws.Range("J3").Value = 30

# Row 4
$ws.Range("D4").Value = 44280

# Row 5
$ws.Range("D5").Value = 44277
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("P5").Value = 550
